$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 01:22"

# Apply updated country stats / re-sorted rows
$ws.Range("B4").Value = 559968
$ws.Range("C4").Value = 27089
$ws.Range("E4").Value = 505946
$ws.Range("G4").Value = 1459
$ws.Range("H4").Value = 22036
$ws.Range("B96").Value = 480
$ws.Range("D96").Value = 231
$ws.Range("E96").Value = 242
$ws.Range("B111").Value = 262
$ws.Range("C111").Value = 4
$ws.Range("E111").Value = 118
$ws.Range("A143").Value = "Bermudas"
$ws.Range("B143").Value = 57
$ws.Range("C143").Value = 9
$ws.Range("D143").Value = 29
$ws.Range("E143").Value = 24
$ws.Range("F143").Value = 2
$ws.Range("H143").Value = 4
$ws.Range("A144").Value = "Uganda"
$ws.Range("B144").Value = 54
$ws.Range("C144").Value = 1
$ws.Range("D144").Value = 4
$ws.Range("E144").Value = 50
$ws.Range("F144").Value = 0
$ws.Range("A145").Value = "Polinesia Francesa"
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 0
$ws.Range("E145").Value = 53
$ws.Range("F145").Value = 1
$ws.Range("H145").Value = 0
$ws.Range("A146").Value = "Islas Caimanes"
$ws.Range("B146").Value = 53
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 6
$ws.Range("E146").Value = 46
$ws.Range("F146").Value = 3
$ws.Range("H146").Value = 1
$ws.Range("A147").Value = "Liberia"
$ws.Range("C147").Value = 2
$ws.Range("D147").Value = 3
$ws.Range("E147").Value = 42
$ws.Range("F147").Value = 0
$ws.Range("H147").Value = 5
$ws.Range("A148").Value = "San Martin (Parte Holandesa)"
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 5
$ws.Range("E148").Value = 36
$ws.Range("H148").Value = 9
$ws.Range("A166").Value = "Mozambique"
$ws.Range("C166").Value = 1
$ws.Range("D166").Value = 2
$ws.Range("F166").Value = 0
$ws.Range("H166").Value = 0
$ws.Range("A167").Value = "Antigua y Barbuda"
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 0
$ws.Range("F167").Value = 1
$ws.Range("H167").Value = 2
$ws.Range("A189").Value = "Nepal"
$ws.Range("C189").Value = 3
$ws.Range("A190").Value = "San Vicente y las Granadinas"
$ws.Range("C190").Value = 0
$ws.Range("A208").Value = "Sudan del Sur"
$ws.Range("A209").Value = "Santo Tome y Principe"
$ws.Range("A210").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("C210").Value = 1
$ws.Range("A211").Value = "Anguila"
$ws.Range("C211").Value = 0
